$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Timesheet")

# Move the timesheet to week 31 (was week 27) - all date formulas recalc automatically.
$ws.Range("C1").Value = 31

# The L8 "comment" dropdown cell now references a new free-text comment "tt"
# instead of the canned "Test execution..." string.
$ws.Range("L8").Value = "tt"

# Agregate Codes Hours: add merged totals, but only for the hours that fall
# within the current month (rows 7-8 => one code, rows 9 and 10 => two more
# codes), summing just the week columns (E:I / E:H) for each.
$ws.Range("C37").Formula = "=SUM(E7:I8)"
$ws.Range("C38").Formula = "=SUM(E9:H9)"
$ws.Range("C39").Formula = "=SUM(E10:H10)"

# Restore the author's last selection/active cell.
$ws.Range("E23").Select()
